$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 92.9375
$ws.Range("I4").Value = 43.444443
$ws.Range("J4").Value = 156.57143
$ws.Range("K4").Value = 43.444443
$ws.Range("L4").Value = 156.57143
$ws.Range("M4").Value = 70.55555699999999
$ws.Range("N4").Value = -384.57143
$ws.Range("H17").Value = 1093.4117
$ws.Range("J17").Value = 1104.36
$ws.Range("L17").Value = 3313.08
$ws.Range("N17").Value = -3649.08
$ws.Range("H19").Value = 1161.5834
$ws.Range("J19").Value = 1571.2858
$ws.Range("L19").Value = 1571.2858
$ws.Range("N19").Value = -1921.2858
$ws.Range("H40").Value = 5256.154
$ws.Range("J40").Value = 5484.8486
$ws.Range("L40").Value = 5484.8486
$ws.Range("N40").Value = -5834.8486
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H106").Value = 3855.7778
$ws.Range("I106").Value = 3855.7778
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3855.7778
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3224.7778
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 3999.6428
$ws.Range("I113").Value = 3713.7144
$ws.Range("K113").Value = 3713.7144
$ws.Range("M113").Value = -459.7143999999998
$ws.Range("H115").Value = 336.8
$ws.Range("I115").Value = 336.8
$ws.Range("K115").Value = 1010.4
$ws.Range("M115").Value = 556.5999999999999
$ws.Range("H116").Value = 12122.833
$ws.Range("I116").Value = 13913.223
$ws.Range("J116").Value = 10332.444
$ws.Range("K116").Value = 13913.223
$ws.Range("L116").Value = 10332.444
$ws.Range("M116").Value = -10471.223
$ws.Range("N116").Value = -17216.444
$ws.Range("H125").Value = 4998.143
$ws.Range("I125").Value = 4749.5
$ws.Range("K125").Value = 42745.5
$ws.Range("M125").Value = -40285.5
$ws.Range("H132").Value = 1322.5555
$ws.Range("I132").Value = 1112.125
$ws.Range("K132").Value = 3336.375
$ws.Range("M132").Value = -806.375
$ws.Range("H135").Value = 6419.091
$ws.Range("J135").Value = 10870.75
$ws.Range("L135").Value = 97836.75
$ws.Range("N135").Value = -102906.75
$ws.Range("H137").Value = 7797.222
$ws.Range("I137").Value = 2417.2
$ws.Range("J137").Value = 34697.332
$ws.Range("K137").Value = 7251.599999999999
$ws.Range("L137").Value = 104091.996
$ws.Range("M137").Value = -4701.599999999999
$ws.Range("N137").Value = -109191.996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1321.5
$ws.Range("I2").Value = 1321.5
$ws.Range("K2").Value = 1321.5
$ws.Range("M2").Value = -1208.5
$ws.Range("H32").Value = 6854.125
$ws.Range("I32").Value = 4288.1943
$ws.Range("K32").Value = 4288.1943
$ws.Range("M32").Value = -4001.1943
$ws.Range("H74").Value = 7444.5527
$ws.Range("J74").Value = 11216.75
$ws.Range("L74").Value = 11216.75
$ws.Range("N74").Value = -12964.75
$ws.Range("H77").Value = 7444.5527
$ws.Range("J77").Value = 11216.75
$ws.Range("L77").Value = 56083.75
$ws.Range("N77").Value = -64819.75
$ws.Range("H116").Value = 1321.5
$ws.Range("I116").Value = 1321.5
$ws.Range("K116").Value = 1321.5
$ws.Range("M116").Value = 972.5
$ws.Range("H132").Value = 5660.5347
$ws.Range("I132").Value = 3800.375
$ws.Range("K132").Value = 11401.125
$ws.Range("M132").Value = -8871.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1321.5
$ws.Range("I3").Value = 1321.5
$ws.Range("K3").Value = 1321.5
$ws.Range("M3").Value = -1207.5
$ws.Range("H20").Value = 41122.43
$ws.Range("I20").Value = 59652.74
$ws.Range("J20").Value = 2002.8889
$ws.Range("K20").Value = 59652.74
$ws.Range("L20").Value = 2002.8889
$ws.Range("M20").Value = -59405.74
$ws.Range("N20").Value = -2496.8889
$ws.Range("H80").Value = 460.94116
$ws.Range("I80").Value = 888.8570999999999
$ws.Range("J80").Value = 161.4
$ws.Range("K80").Value = 888.8570999999999
$ws.Range("L80").Value = 161.4
$ws.Range("M80").Value = 109.1429000000001
$ws.Range("N80").Value = -2157.4
$ws.Range("H83").Value = 460.94116
$ws.Range("I83").Value = 888.8570999999999
$ws.Range("J83").Value = 161.4
$ws.Range("K83").Value = 4444.2855
$ws.Range("L83").Value = 807
$ws.Range("M83").Value = 547.7145
$ws.Range("N83").Value = -10791
$ws.Range("H134").Value = 1111.826
$ws.Range("I134").Value = 943
$ws.Range("K134").Value = 2829
$ws.Range("M134").Value = -294
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2115.375
$ws.Range("I22").Value = 780.82355
$ws.Range("K22").Value = 780.82355
$ws.Range("M22").Value = -430.82355
$ws.Range("H122").Value = 19559.928
$ws.Range("I122").Value = 1512.6364
$ws.Range("J122").Value = 85733.336
$ws.Range("K122").Value = 4537.9092
$ws.Range("L122").Value = 257200.008
$ws.Range("M122").Value = -2087.9092
$ws.Range("N122").Value = -262100.008
$ws.Range("H132").Value = 6426.517
$ws.Range("I132").Value = 7414.591
$ws.Range("K132").Value = 22243.773
$ws.Range("M132").Value = -19713.773
$ws.Range("H134").Value = 4959
$ws.Range("I134").Value = 4888.5
$ws.Range("K134").Value = 14665.5
$ws.Range("M134").Value = -12130.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 502551.9
$ws.Range("J129").Value = 2605.3845
$ws.Range("L129").Value = 7816.1535
$ws.Range("N129").Value = -17816.1535
$ws.Range("H132").Value = 5607.8887
$ws.Range("J132").Value = 6475.2383
$ws.Range("L132").Value = 58277.1447
$ws.Range("N132").Value = -63337.1447
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 13001
$ws.Range("I10").Value = 20000
$ws.Range("J10").Value = 9501.5
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 9501.5
$ws.Range("M10").Value = -19831
$ws.Range("N10").Value = -9839.5
$ws.Range("H48").Value = 5640000.5
$ws.Range("J48").Value = 30001
$ws.Range("L48").Value = 30001
$ws.Range("N48").Value = -30971
$ws.Range("H102").Value = 2517.7058
$ws.Range("I102").Value = 1020.06665
$ws.Range("K102").Value = 1020.06665
$ws.Range("M102").Value = 601.93335
$ws.Range("H126").Value = 3289.5881
$ws.Range("I126").Value = 3372.5
$ws.Range("J126").Value = 3171.1428
$ws.Range("K126").Value = 10117.5
$ws.Range("L126").Value = 9513.428400000001
$ws.Range("M126").Value = -7647.5
$ws.Range("N126").Value = -14453.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 52000
$ws.Range("I74").Value = 52000
$ws.Range("K74").Value = 52000
$ws.Range("M74").Value = -51002
$ws.Range("H76").Value = 18499
$ws.Range("I76").Value = 9998
$ws.Range("K76").Value = 9998
$ws.Range("M76").Value = -9660
$ws.Range("H77").Value = 52000
$ws.Range("I77").Value = 52000
$ws.Range("K77").Value = 156000
$ws.Range("M77").Value = -151008
$ws.Range("H79").Value = 18499
$ws.Range("I79").Value = 9998
$ws.Range("K79").Value = 9998
$ws.Range("M79").Value = -8828
$ws.Range("H93").Value = 9433.286
$ws.Range("I93").Value = 8516
$ws.Range("J93").Value = 10656.333
$ws.Range("K93").Value = 8516
$ws.Range("L93").Value = 10656.333
$ws.Range("M93").Value = -7268
$ws.Range("N93").Value = -13152.333
$ws.Range("H106").Value = 12253.8
$ws.Range("J106").Value = 12253.8
$ws.Range("L106").Value = 12253.8
$ws.Range("N106").Value = -14777.8
$ws.Range("H132").Value = 4523.75
$ws.Range("I132").Value = 4570
$ws.Range("K132").Value = 13710
$ws.Range("M132").Value = -11180
$ws.Range("H136").Value = 7105.5454
$ws.Range("I136").Value = 7116.1
$ws.Range("K136").Value = 21348.3
$ws.Range("M136").Value = -18798.3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1850.725
$ws.Range("I113").Value = 1378.4138
$ws.Range("K113").Value = 4135.2414
$ws.Range("M113").Value = -1965.2414
$ws.Range("H122").Value = 51258.39
$ws.Range("I122").Value = 1857.2142
$ws.Range("J122").Value = 128104.664
$ws.Range("K122").Value = 5571.642599999999
$ws.Range("L122").Value = 384313.992
$ws.Range("M122").Value = -3121.642599999999
$ws.Range("N122").Value = -389213.992
$ws.Range("H132").Value = 1768.5135
$ws.Range("I132").Value = 1678.7778
$ws.Range("K132").Value = 5036.3334
$ws.Range("M132").Value = -2506.3334
$ws.Range("H133").Value = 90325
$ws.Range("I133").Value = 80650
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 80650
$ws.Range("L133").Value = 100000
$ws.Range("M133").Value = -75590
$ws.Range("N133").Value = -110120
$ws.Range("H136").Value = 69512.87
$ws.Range("I136").Value = 2485.3635
$ws.Range("K136").Value = 7456.0905
$ws.Range("M136").Value = -4906.0905
